$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '30.414.62'
Set-TextValue "D3" '1.917.41'
Set-TextValue "E3" '  +2.02%  '
Set-TextValue "E4" '  -0.15%  '
Set-TextValue "D5" '241.07'
Set-TextValue "E5" '  +1.51%  '
Set-TextValue "D6" '1.000'
Set-TextValue "E6" '  -0.13%  '
Set-TextValue "D7" '0.4689'
Set-TextValue "E7" '  -1.40%  '
Set-TextValue "E8" '  +0.79%  '
Set-TextValue "D9" '0.06809'
Set-TextValue "E9" '  +4.87%  '
Set-TextValue "D10" '107.26'
Set-TextValue "E11" '  -2.57%  '
Set-TextValue "D12" '1.901.40'
Set-TextValue "E12" '  +1.09%  '
Set-TextValue "E13" '  +0.67%  '
Set-TextValue "D14" '5.187'
Set-TextValue "E14" '  +2.97%  '
Set-TextValue "D15" '0.6543'
Set-TextValue "E15" '  +1.05%  '
Set-TextValue "D16" '288.85'
Set-TextValue "E16" '  -2.96%  '
Set-TextValue "D17" '30.402.65'
Set-TextValue "E17" '  -0.91%  '
Set-TextValue "B18" 'ShibaInu'
Set-TextValue "C18" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D18" '0.000007628'
Set-TextValue "E18" '  +1.85%  '
Set-TextValue "B19" 'Avalanche'
Set-TextValue "C19" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D19" '12.96'
Set-TextValue "E19" '  -0.60%  '
Set-TextValue "D20" '0.9994'
Set-TextValue "E20" '  -0.24%  '
Set-TextValue "D21" '2.150.10'
Set-TextValue "E21" '  +1.03%  '
Set-TextValue "D22" '1.001'
Set-TextValue "E22" '  +0.04%  '
Set-TextValue "D23" '5.218'
Set-TextValue "E23" '  +1.89%  '
Set-TextValue "D24" '6.198'
Set-TextValue "E24" '  +1.51%  '
Set-TextValue "D25" '21.72'
Set-TextValue "E25" '  +11.35%  '
Set-TextValue "D26" '168.04'
Set-TextValue "E26" '  -0.72%  '
Set-TextValue "D27" '9.262'
Set-TextValue "E27" '  +0.58%  '
Set-TextValue "D28" '2.040'
Set-TextValue "E28" '  +4.97%  '
Set-TextValue "D29" '0.1069'
Set-TextValue "E29" '  +1.24%  '
Set-TextValue "D30" '1.371'
Set-TextValue "E30" '  +1.64%  '
Set-TextValue "D31" '4.141'
Set-TextValue "E31" '  +0.01%  '
Set-TextValue "D32" '3.939'
Set-TextValue "E32" '  +0.02%  '
Set-TextValue "D33" '0.05026'
Set-TextValue "E33" '  +1.17%  '
Set-TextValue "D34" '0.7376'
Set-TextValue "E34" '  +3.01%  '
Set-TextValue "D35" '1.147'
Set-TextValue "E35" '  -1.59%  '
Set-TextValue "D36" '2.739'
Set-TextValue "E36" '  +0.79%  '
Set-TextValue "D37" '0.02031'
Set-TextValue "E37" '  +6.81%  '
Set-TextValue "E38" '  -0.50%  '
Set-TextValue "D39" '2.049'
Set-TextValue "E39" '  +0.28%  '
Set-TextValue "D40" '108.56'
Set-TextValue "E40" '  +1.66%  '
Set-TextValue "D41" '0.8729'
Set-TextValue "E41" '  -2.53%  '
Set-TextValue "D42" '5.846'
Set-TextValue "E42" '  +5.05%  '
Set-TextValue "D43" '1.0000'
Set-TextValue "E43" '  -0.14%  '
Set-TextValue "D44" '52.97'
Set-TextValue "E44" '  +26.59%  '
Set-TextValue "D45" '0.4204'
Set-TextValue "E45" '  +0.77%  '
Set-TextValue "D46" '67.51'
Set-TextValue "E46" '  +3.35%  '
Set-TextValue "E47" '  -2.00%  '
Set-TextValue "D48" '9.215'
Set-TextValue "E48" '  +4.81%  '
Set-TextValue "D49" '0.1206'
Set-TextValue "E49" '  -0.57%  '
Set-TextValue "D50" '34.67'
Set-TextValue "E50" '  +0.71%  '
Set-TextValue "D51" '0.00004407'
Set-TextValue "E51" '  +40.40%  '
